# "addback some models, create others"
# - Options sheet: lower the "PRICE ON FAILURE" assumption (H10) from 3 to 2,
#   and add back a "12/20/24 2.5 PUTS" row above the existing 12/20/24 put
#   ladder (also tweaking a couple of the existing strike costs).
# - View-state: Main's selection moves off the old A1:K13 block, FLNA's
#   selection resets to B1, and the Options sheet becomes the active tab.

$wb = $excel.ActiveWorkbook
$wsOptions = $wb.Worksheets.Item("Options")

# Lower the "PRICE ON FAILURE" input.
$wsOptions.Range("H10").Value = 2

# Insert a new row above the first 12/20/24 PUTS row (row 13), pushing the
# rest of the ladder down by one and inheriting the header row's formats.
[void]$wsOptions.Rows("13:13").Insert()

$wsOptions.Range("B13").Value = "12/20/24 2.5 PUTS"
$wsOptions.Range("C13").Value = 0.25
$wsOptions.Range("E13").Value = 2.5
$wsOptions.Range("F13").Formula = "=+E13-`$H`$10"
$wsOptions.Range("G13").Formula = "=+F13/C13-1"

# Match the number formats/fonts of the row below (the template row Excel
# copied down from the header lost the ladder's numeric formatting).
[void]$wsOptions.Range("F14").Copy()
[void]$wsOptions.Range("F13").PasteSpecial(-4122)
[void]$wsOptions.Range("G14").Copy()
[void]$wsOptions.Range("G13").PasteSpecial(-4122)

# Update the costs that were re-priced on the (now shifted) existing rows.
$wsOptions.Range("C15").Value = 1.75
$wsOptions.Range("C16").Value = 2.9

# --- view-state / selection updates ---
$wsMain = $wb.Worksheets.Item("Main")
[void]$wsMain.Range("J5").Select()

$wsFLNA = $wb.Worksheets.Item("FLNA")
[void]$wsFLNA.Range("B1").Select()

[void]$wsOptions.Activate()
[void]$wsOptions.Range("A16").Select()
